# Commit: "Fruta / hortaliza, semanal"
# A new weekly price record (Naranja / Navel Late / Primera, 2021-12-29) is inserted
# into the data table at row 293, pushing the existing rows 293:387 down to 294:388.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 293; this shifts rows 293-387 down
# to 294-388 and automatically extends the sheet's used range / dimension.
$ws.Rows.Item(293).Insert()

# Populate the newly-inserted row 293 with the new record.
$ws.Range("A293").Value = 2
$ws.Range("B293").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C293").Value = "Coquimbo"
$ws.Range("D293").Value = 44559
$ws.Range("E293").Value = 4
$ws.Range("F293").Value = "Fruta"
$ws.Range("G293").Value = 100102
$ws.Range("H293").Value = "Cítricos"
$ws.Range("I293").Value = 100102005
$ws.Range("J293").Value = "Naranja"
$ws.Range("K293").Value = "Navel Late"
$ws.Range("L293").Value = "Primera"
$ws.Range("M293").Value = 24
$ws.Range("N293").Value = 225000
$ws.Range("O293").Value = 230000
$ws.Range("P293").Value = 227500
$ws.Range("Q293").Value = "$/bins (400 kilos)"
$ws.Range("R293").Value = "Provincia de Limarí"
$ws.Range("S293").Value = 569
$ws.Range("T293").Value = 400
